# Weekly refresh of "Fruta, Vega Central Mapocho de Santiago - Coco" data:
# the per-record fields (Fecha, Calidad, Volumen, Precio minimo/maximo/promedio,
# Origen, Precio $/Kg) are reshuffled across the existing data rows (2..39).
# Columns A,B,C,E,F,G,H,I,J,K,Q,T are identical for every row already, so the
# only observable effect is columns D,L,M,N,O,P,R,S moving to a different row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row (2..39) -> source row whose D/L/M/N/O/P/R/S values it now holds
$srcForTarget = @{
    2  = 10
    3  = 6
    4  = 26
    5  = 24
    6  = 38
    7  = 9
    8  = 3
    9  = 25
    10 = 15
    11 = 20
    12 = 32
    13 = 23
    14 = 28
    15 = 33
    16 = 37
    17 = 12
    18 = 17
    19 = 29
    20 = 8
    21 = 34
    22 = 19
    23 = 14
    24 = 18
    25 = 7
    26 = 13
    27 = 27
    28 = 21
    29 = 11
    30 = 30
    31 = 39
    32 = 22
    33 = 16
    34 = 5
    35 = 31
    36 = 4
    37 = 35
    38 = 2
    39 = 36
}

$cols = @(4, 12, 13, 14, 15, 16, 18, 19)   # D, L, M, N, O, P, R, S

# Snapshot the original values of the moving columns for every data row
# before writing anything, since the permutation has long cycles and
# source rows must still hold their original values when they are read.
$snapshot = @{}
for ($r = 2; $r -le 39; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

for ($r = 2; $r -le 39; $r++) {
    $src = $srcForTarget[$r]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $snapshot[$src][$c]
    }
}
